$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The EC (Estado de Cuenta) database is updated: the two workers that were in
# rows 16 and 18 swap places (row 17 - SOL AISEA VILLA MARTINEZ - is untouched),
# and their overdue/base salary figures move together with them.
#
# Row 16 -> now YENIS PAOLA PALACIO CASTELLAR (CC 1143325780), periodo 1702,
#           Valor Mora 4961, Salario Basico 1240200
# Row 18 -> now DIANA MARCELA PARRA MORALES (CC 1143373551), periodo 2405,
#           Valor Mora 88000, Salario Basico 2200000

$ws.Range("C16").Value = "1143325780"
$ws.Range("D16").Value = "YENIS PAOLA PALACIO CASTELLAR"
$ws.Range("E16").Value = "1702"
$ws.Range("F16").Value = 4961
$ws.Range("G16").Value = 1240200

$ws.Range("C18").Value = "1143373551"
$ws.Range("D18").Value = "DIANA MARCELA PARRA MORALES"
$ws.Range("E18").Value = "2405"
$ws.Range("F18").Value = 88000
$ws.Range("G18").Value = 2200000

$wb.Save()
